# Replace the numeric month values in column C (rows 5-84) with their
# Spanish three-letter month-abbreviation text equivalents (e.g. 8 -> "Ago.")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_2 y 3")

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 5; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)   # column C
    $monthNum = [int]$cell.Value()
    $cell.Value = $monthNames[$monthNum]
}
